$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 782
$ws.Cells.Item(2, 2).Value = "Chronic Thromboembolic Pulmonary Hypertension"
